$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.271.88"
$ws.Range("E2").Value = "  -1.10%  "
$ws.Range("D3").Value = "2.247.49"
$ws.Range("E3").Value = "  -1.32%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.19"
$ws.Range("E5").Value = "  -1.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.620"
$ws.Range("E6").Value = "  -3.52%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.47"
$ws.Range("E7").Value = "  -0.92%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.620"
$ws.Range("E9").Value = "  -2.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.10"
$ws.Range("E10").Value = "  +7.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0944"
$ws.Range("E11").Value = "  -3.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.15"
$ws.Range("E12").Value = "  -3.58%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.103"
$ws.Range("E13").Value = "  -2.65%  "
$ws.Range("D14").Value = "2.582.94"
$ws.Range("E14").Value = "  -1.36%  "
$ws.Range("E15").Value = "  -3.46%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.856"
$ws.Range("E16").Value = "  -1.22%  "
$ws.Range("D17").Value = "2.255.08"
$ws.Range("E17").Value = "  -0.85%  "
$ws.Range("D18").Value = "42.123.29"
$ws.Range("E18").Value = "  -1.21%  "
$ws.Range("D19").Value = "0.0₃0980"
$ws.Range("E19").Value = "  -1.72%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.11"
$ws.Range("E20").Value = "  -1.53%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.91"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.26"
$ws.Range("E22").Value = "  +4.43%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "231.84"
$ws.Range("E23").Value = "  -1.20%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.30"
$ws.Range("E25").Value = "  +30.88%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.18"
$ws.Range("E26").Value = "  -0.64%  "
$ws.Range("E27").Value = "  -7.63%  "
$ws.Range("E28").Value = "  -3.51%  "
$ws.Range("E29").Value = "  +1.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "169.37"
$ws.Range("E30").Value = "  +1.13%  "
$ws.Range("E31").Value = "  -1.63%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0823"
$ws.Range("E32").Value = "  -5.72%  "
$ws.Range("E33").Value = "  -5.79%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "30.46"
$ws.Range("E34").Value = "  -4.34%  "
$ws.Range("E35").Value = "  -1.93%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.56"
$ws.Range("E36").Value = "  +0.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.98"
$ws.Range("E37").Value = "  +4.76%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0307"
$ws.Range("E38").Value = "  +0.80%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.48"
$ws.Range("E39").Value = "  -0.50%  "
$ws.Range("E40").Value = "  -4.60%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.79"
$ws.Range("E41").Value = "  -1.37%  "
$ws.Range("E42").Value = "  -2.57%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "61.32"
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "107.82"
$ws.Range("E44").Value = "  +1.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.64"
$ws.Range("E45").Value = "  -3.08%  "
$ws.Range("E46").Value = "  -0.11%  "
$ws.Range("E47").Value = "  -0.36%  "
$ws.Range("E48").Value = "  -3.76%  "
$ws.Range("E49").Value = "  -0.94%  "
$ws.Range("E50").Value = "  +1.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.10"
$ws.Range("E51").Value = "  -2.55%  "
